$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19, shifting existing rows 19-45 down to 20-46.
$ws.Rows(19).Insert()

# New row 19 content: NATLIFEINS / 3. Copy the cell format from the row below
# (now row 20, formerly row 19) so it matches the other ticker cells in column A.
$ws.Range("A20").Copy()
$ws.Range("A19").PasteSpecial(-4122)
$ws.Range("A19").Value = "NATLIFEINS"
$ws.Range("B19").Value = 3

# Update dimension-affecting weight values to match final state.
$ws.Range("B2").Value = 2
$ws.Range("B8").Value = 6
$ws.Range("B9").Value = 1
$ws.Range("B10").Value = 5
$ws.Range("B11").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B21").Value = 4
$ws.Range("B23").Value = 2
$ws.Range("B24").Value = 3
$ws.Range("B25").Value = 3
$ws.Range("B26").Value = 1
$ws.Range("B29").Value = 1
$ws.Range("B30").Value = 1
$ws.Range("B31").Value = 2
$ws.Range("B34").Value = 2
$ws.Range("B35").Value = 3
$ws.Range("B38").Value = 3
$ws.Range("B39").Value = 1.5
$ws.Range("B40").Value = 1.5
$ws.Range("B41").Value = 3
$ws.Range("B42").Value = 3
$ws.Range("B45").Value = 3
$ws.Range("B46").Value = 2
